$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.052.55"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "1.888.19"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7396"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.67"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3183"
$ws.Range("E8").Value = "  +2.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07193"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.80"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08309"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7610"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.408"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").Value = "1.869.49"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.15"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "30.071.78"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.168"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "252.52"
$ws.Range("E18").Value = "  +5.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.59"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007918"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.130.99"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9990"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.873"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1551"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.283"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.45"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.75"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.051"
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.465"
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.586"
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.535"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.217"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.252"
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7679"
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9983"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.737"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.762"
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4587"
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.044"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.56"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8723"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.077.69"
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.38"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.861"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.051.88"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.897"
$ws.Range("E51").Value = "  -3.80%  "
